$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Jun2020")

$data = @(
    @(2, "Lem", "Momay"),
    @(3, "Nichole", "Nathan"),
    @(4, "Ardeth", "Anj"),
    @(5, "Kate", "Lem"),
    @(6, "Camille", "Dan"),
    @(7, "Kevin", "John"),
    @(8, "Carlo", "Amy"),
    @(9, "Oscar", "Reyn"),
    @(10, "Kevin", "Nora"),
    @(11, "Roy", "Jec"),
    @(12, "Sky", "Brain"),
    @(13, "Morgan", "Web"),
    @(14, "Tine", "Carlo"),
    @(15, "Momay", "Harry"),
    @(16, "Troy", "Arnel"),
    @(17, "Issa", "Eunice"),
    @(18, "Anne", "Cath"),
    @(19, "Gene", "Kate"),
    @(20, "Eunice", "Louie"),
    @(21, "Nora", "Drew"),
    @(22, "Rodney", "Troy"),
    @(23, "Meryll", "David"),
    @(24, "Arnel", "Nora"),
    @(25, "Janine", "Tina"),
    @(26, "Cath", "Sky"),
    @(27, "Web", "Anne"),
    @(28, "Julius", "Dan"),
    @(29, "Anj", "Gene"),
    @(30, "James", "Johann"),
    @(31, "John", "Meryll"),
    @(32, "Drew", "Daisy"),
    @(33, "Carla", "Rodney"),
    @(34, "Chad", "Pau"),
    @(35, "Luz", "Nichole"),
    @(36, "Nathan", "Camille"),
    @(37, "David", "Roy"),
    @(38, "Tina", "Oscar"),
    @(39, "Harry", "Janine"),
    @(40, "Pau", "Web"),
    @(41, "John", "Kevin"),
    @(42, "Luz", "Jec"),
    @(43, "Louie", "Chad"),
    @(44, "Harle", "James"),
    @(45, "Nathan", "Kennex"),
    @(46, "Ken", "Carla"),
    @(47, "Jes", "Ardeth"),
    @(48, "Nichole", "Luz"),
    @(49, "JK", "Eunice"),
    @(50, "Jes", "Tintin"),
    @(51, "Harle", "Morgan"),
    @(52, "Cath", "Julius"),
    @(53, "Johann", "Jes"),
    @(54, "Sky", "Issa"),
    @(55, "Jec", "Tine"),
    @(56, "Reyn", "Harry"),
    @(57, "Kennex", "Harle"),
    @(58, "Dan", "Jhoanne"),
    @(59, "Brain", "Ken"),
    @(60, "James", "Reyn"),
    @(61, "Gene", "JK")
)

foreach ($row in $data) {
    $r = $row[0]
    $ws.Cells.Item($r, 2).Value = $row[1]
    $ws.Cells.Item($r, 3).Value = $row[2]
}

$ws.Activate()
$ws.Range("C8").Select()
